$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.420.70"
Set-TextValue $ws.Range("E2") "  +0.02%  "
Set-TextValue $ws.Range("D3") "1.854.01"
Set-TextValue $ws.Range("E3") "  +0.33%  "
Set-TextValue $ws.Range("D4") "0.9995"
Set-TextValue $ws.Range("E4") "  -0.01%  "
Set-TextValue $ws.Range("D5") "241.21"
Set-TextValue $ws.Range("E5") "  +0.20%  "
Set-TextValue $ws.Range("D6") "0.6342"
Set-TextValue $ws.Range("E6") "  +0.76%  "
Set-TextValue $ws.Range("E7") "  +0.04%  "
Set-TextValue $ws.Range("D8") "0.07589"
Set-TextValue $ws.Range("E8") "  -1.58%  "
Set-TextValue $ws.Range("D9") "0.2931"
Set-TextValue $ws.Range("E9") "  +0.31%  "
Set-TextValue $ws.Range("D10") "24.55"
Set-TextValue $ws.Range("E10") "  -1.02%  "
Set-TextValue $ws.Range("D11") "0.07756"
Set-TextValue $ws.Range("E11") "  +0.21%  "
Set-TextValue $ws.Range("D12") "1.853.18"
Set-TextValue $ws.Range("E12") "  +0.32%  "
Set-TextValue $ws.Range("D13") "5.030"
Set-TextValue $ws.Range("E13") "  +0.02%  "
Set-TextValue $ws.Range("D14") "0.6822"
Set-TextValue $ws.Range("E14") "  +0.28%  "
Set-TextValue $ws.Range("D15") "0.00001050"
Set-TextValue $ws.Range("E15") "  -2.06%  "
Set-TextValue $ws.Range("D16") "83.37"
Set-TextValue $ws.Range("E16") "  -0.20%  "
Set-TextValue $ws.Range("B17") "WrappedliquidstakedEther2.0"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D17") "2.107.74"
Set-TextValue $ws.Range("E17") "  +0.43%  "
Set-TextValue $ws.Range("B18") "Uniswap"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D18") "6.149"
Set-TextValue $ws.Range("E18") "  -0.46%  "
Set-TextValue $ws.Range("B19") "WrappedBTC"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D19") "29.433.62"
Set-TextValue $ws.Range("E19") "  -0.06%  "
Set-TextValue $ws.Range("B20") "BitcoinCash"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D20") "230.55"
Set-TextValue $ws.Range("E20") "  +1.02%  "
Set-TextValue $ws.Range("B21") "Avalanche"
Set-TextValue $ws.Range("C21") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D21") "12.39"
Set-TextValue $ws.Range("E21") "  -0.26%  "
Set-TextValue $ws.Range("B22") "Dai"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D22") "1.000"
Set-TextValue $ws.Range("E22") "  -0.04%  "
Set-TextValue $ws.Range("B23") "Chainlink"
Set-TextValue $ws.Range("C23") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D23") "7.496"
Set-TextValue $ws.Range("E23") "  +1.03%  "
Set-TextValue $ws.Range("B24") "BinanceUSD"
Set-TextValue $ws.Range("C24") "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D24") "1.000"
Set-TextValue $ws.Range("E24") "  -0.04%  "
Set-TextValue $ws.Range("B25") "Monero"
Set-TextValue $ws.Range("C25") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D25") "159.23"
Set-TextValue $ws.Range("E25") "  +1.05%  "
Set-TextValue $ws.Range("B26") "Stellar"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D26") "0.1399"
Set-TextValue $ws.Range("E26") "  +1.35%  "
Set-TextValue $ws.Range("B27") "Cosmos"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D27") "8.466"
Set-TextValue $ws.Range("E27") "  +0.73%  "
Set-TextValue $ws.Range("B28") "EthereumClassic"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D28") "17.70"
Set-TextValue $ws.Range("E28") "  +0.10%  "
Set-TextValue $ws.Range("B29") "Toncoin"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D29") "1.417"
Set-TextValue $ws.Range("E29") "  +4.46%  "
Set-TextValue $ws.Range("B30") "PancakeSwap"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D30") "1.477"
Set-TextValue $ws.Range("E30") "  +0.98%  "
Set-TextValue $ws.Range("B31") "Hedera"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D31") "0.05700"
Set-TextValue $ws.Range("E31") "  +0.22%  "
Set-TextValue $ws.Range("B32") "Filecoin"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D32") "4.133"
Set-TextValue $ws.Range("E32") "  +0.34%  "
Set-TextValue $ws.Range("B33") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D33") "4.057"
Set-TextValue $ws.Range("E33") "  +0.69%  "
Set-TextValue $ws.Range("B34") "LidoDAOToken"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D34") "1.832"
Set-TextValue $ws.Range("E34") "  -0.56%  "
Set-TextValue $ws.Range("B35") "ARBITRUM"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D35") "1.158"
Set-TextValue $ws.Range("E35") "  -0.24%  "
Set-TextValue $ws.Range("B36") "ImmutableX"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D36") "0.6992"
Set-TextValue $ws.Range("E36") "  -1.25%  "
Set-TextValue $ws.Range("B37") "HuobiToken"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D37") "2.584"
Set-TextValue $ws.Range("E37") "  +0.02%  "
Set-TextValue $ws.Range("B38") "VeChain"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D38") "0.01833"
Set-TextValue $ws.Range("E38") "  +2.43%  "
Set-TextValue $ws.Range("B39") "Maker"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D39") "1.248.36"
Set-TextValue $ws.Range("E39") "  +2.02%  "
Set-TextValue $ws.Range("B40") "MXToken"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D40") "2.728"
Set-TextValue $ws.Range("E40") "  -1.84%  "
Set-TextValue $ws.Range("B41") "FraxShare"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D41") "6.462"
Set-TextValue $ws.Range("E41") "  -1.24%  "
Set-TextValue $ws.Range("B42") "TrustWalletToken"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D42") "0.9039"
Set-TextValue $ws.Range("E42") "  -0.35%  "
Set-TextValue $ws.Range("B43") "PaxDollar"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D43") "1.000"
Set-TextValue $ws.Range("E43") "  -0.02%  "
Set-TextValue $ws.Range("B44") "RocketPoolETH"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue $ws.Range("D44") "2.012.18"
Set-TextValue $ws.Range("E44") "  +0.22%  "
Set-TextValue $ws.Range("B45") "Quant"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D45") "101.94"
Set-TextValue $ws.Range("E45") "  +0.17%  "
Set-TextValue $ws.Range("B46") "Aave"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D46") "66.05"
Set-TextValue $ws.Range("E46") "  -0.24%  "
Set-TextValue $ws.Range("B47") "Aptos"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D47") "7.155"
Set-TextValue $ws.Range("E47") "  -0.04%  "
Set-TextValue $ws.Range("B48") "Algorand"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D48") "0.1170"
Set-TextValue $ws.Range("E48") "  +2.00%  "
Set-TextValue $ws.Range("D49") "9.017"
Set-TextValue $ws.Range("E49") "  +0.19%  "
Set-TextValue $ws.Range("D50") "0.3972"
Set-TextValue $ws.Range("E50") "  -1.17%  "
Set-TextValue $ws.Range("B51") "BabyDogeCoin"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D51") "0.00000000115"
Set-TextValue $ws.Range("E51") "  -3.76%  "

Write-Host "Applied $($wb.Name) updates"
